$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the old row 806 (the 2026/12/29 block),
# shifting all subsequent rows down by 2. This makes room for two new
# data points for 2026/02/16.
$ws.Rows("806:807").Insert()

# Row 806: 2026/02/16, 月, 19:00, rank 49
$ws.Cells.Item(806, 1).NumberFormat = "@"
$ws.Cells.Item(806, 1).Value = "2026/02/16"
$ws.Cells.Item(806, 1).Style = "Normal"
$ws.Cells.Item(806, 2).NumberFormat = "@"
$ws.Cells.Item(806, 2).Value = "月"
$ws.Cells.Item(806, 2).Style = "Normal"
$ws.Cells.Item(806, 3).Value = 19
$ws.Cells.Item(806, 4).Value = 49

# Row 807: 2026/02/16, 月, 22:00, rank 53
$ws.Cells.Item(807, 1).NumberFormat = "@"
$ws.Cells.Item(807, 1).Value = "2026/02/16"
$ws.Cells.Item(807, 1).Style = "Normal"
$ws.Cells.Item(807, 2).NumberFormat = "@"
$ws.Cells.Item(807, 2).Value = "月"
$ws.Cells.Item(807, 2).Style = "Normal"
$ws.Cells.Item(807, 3).Value = 22
$ws.Cells.Item(807, 4).Value = 53
